$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("constant")

# --- Insert 4 new rows to split progression_rate / containment_rate into age bands ---
# Row 11 currently holds progression_rate; insert two blank rows right after it.
$ws.Rows.Item(12).Insert()
$ws.Rows.Item(12).Insert()

# After the inserts above, the old row 12 (containment_rate) now lives at row 14.
# Insert two more blank rows right after it.
$ws.Rows.Item(15).Insert()
$ws.Rows.Item(15).Insert()

# --- Row 11: progression_rate -> progression_rate_age0 (drop its uniform distribution) ---
$ws.Range("A11").Value = "progression_rate_age0"
$ws.Range("B11").Value = 2.4
$ws.Range("C11:E11").ClearContents()

# --- Row 12 (new): progression_rate_age5 ---
$ws.Range("A12").Value = "progression_rate_age5"
$ws.Range("B12").Value = 2

# --- Row 13 (new): progression_rate_age15 ---
$ws.Range("A13").Value = "progression_rate_age15"
$ws.Range("B13").Value = 0.1

# --- Row 14: containment_rate -> containment_rate_age0 (drop its uniform distribution) ---
$ws.Range("A14").Value = "containment_rate_age0"
$ws.Range("B14").Value = 4.4
$ws.Range("C14:E14").ClearContents()

# --- Row 15 (new): containment_rate_age5 ---
$ws.Range("A15").Value = "containment_rate_age5"
$ws.Range("B15").Value = 4.4

# --- Row 16 (new): containment_rate_age15 ---
$ws.Range("A16").Value = "containment_rate_age15"
$ws.Range("B16").Value = 2

# --- Row 18: clearance_rate now gets a uniform distribution ---
$ws.Range("C18").Value = "uniform"
$ws.Range("D18").Value = 0.01
$ws.Range("E18").Value = 0.03

# --- Row 28 (was row 24): pct_neg_tx_death value changes from 50 to 40 ---
$ws.Range("B28").Value = 40

# --- Update selection to match authored state ---
$ws.Range("B11").Select()
